$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clebina")
$ws.Range("G9").Value = 220
$ws.Range("G10").Value = 6600
